# Add a new player row (Al Horford) and reorder the roster so that
# players are grouped by position, matching the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Austin Reaves",       "PG,SG",    "Los Angeles Lakers"),
    @("Trey Murphy III",     "SG,SF,PF", "New Orleans Pelicans"),
    @("Stephen Curry",       "PG,SG",    "Golden State Warriors"),
    @("Tyrese Haliburton",   "PG,SG",    "Indiana Pacers"),
    @("Cameron Johnson",     "SF,PF",    "Brooklyn Nets"),
    @("Keegan Murray",       "SF,PF",    "Sacramento Kings"),
    @("Franz Wagner",        "SF,PF",    "Orlando Magic"),
    @("Kevin Durant",        "SF,PF",    "Phoenix Suns"),
    @("Al Horford",          "PF,C",     "Boston Celtics"),
    @("Karl-Anthony Towns",  "PF,C",     "New York Knicks"),
    @("Jalen Duren",         "C",        "Detroit Pistons"),
    @("Jarrett Allen",       "C",        "Cleveland Cavaliers"),
    @("Darius Garland",      "PG",       "Cleveland Cavaliers"),
    @("Klay Thompson",       "SG,SF",    "Dallas Mavericks"),
    @("OG Anunoby",          "SF,PF",    "New York Knicks"),
    @("Tyrese Maxey",        "PG,SG",    "Philadelphia 76ers"),
    @("Mark Williams",       "C",        "Charlotte Hornets")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
